$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: clone row 7's formatting/values (same trait "Area", same accession BGV000928 pattern)
# then correct plant_number/value/timestamp to the new observation's data.
$ws.Range("A7:G7").Copy($ws.Range("A12:G12"))
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 2000
$ws.Range("G12").Value = 42528.5661226852

# Row 13: same as row 12 but with a new trait "Growth habit".
$ws.Range("A7:G7").Copy($ws.Range("A13:G13"))
$ws.Range("C13").Value = "Growth habit"
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 2000
$ws.Range("G13").Value = 42528.5661226852

$ws.Range("C13").Select()
